# Swap the data (columns B through AC) between paired rows.
# The "id" column (A) stays associated with its own row; only the
# match-data columns were transposed between the two rows in each pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(383, 384),
    @(408, 409),
    @(410, 411),
    @(420, 421),
    @(441, 443),
    @(453, 454),
    @(464, 465),
    @(469, 470),
    @(476, 477),
    @(505, 506),
    @(571, 572),
    @(575, 576),
    @(739, 740),
    @(742, 743),
    @(772, 773),
    @(775, 777)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AC$r1")
    $range2 = $ws.Range("B$r2`:AC$r2")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}
